$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates (Gdf6-Bmpr1a, Target cluster: ECs)
$ws.Range("M2").Value = 2.341355666666667
$ws.Range("N2").Value = 7.024067000000001
$ws.Range("O2").Value = 0.03973512964576821
$ws.Range("P2").Value = 0.0397351296457682
$ws.Range("Q2").Value = 0.841153095451
$ws.Range("R2").Value = 7.570377859059001
$ws.Range("S2").Value = 0.03973512964576821
$ws.Range("T2").Value = 0.0397351296457682

# Row 3 updates (Gdf6-Bmpr1a, Target cluster: FAPs)
$ws.Range("O3").Value = 0.5779093692199981
$ws.Range("P3").Value = 0.5779093692199981
$ws.Range("S3").Value = 0.5779093692199981
$ws.Range("T3").Value = 0.5779093692199981

# Row 4 updates (Gdf6-Bmpr1a, Target cluster: MuSCs)
$ws.Range("O4").Value = 0.3823555011342337
$ws.Range("P4").Value = 0.3823555011342337
$ws.Range("S4").Value = 0.3823555011342337
$ws.Range("T4").Value = 0.3823555011342337
